$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 171 (pushes existing rows 171-176 down to 172-177).
# The new row duplicates row 170's data but with an updated date (Fecha).
$ws.Rows(171).Insert()

$ws.Cells.Item(171, 1).Value = 2
$ws.Cells.Item(171, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(171, 3).Value = "Coquimbo"
$ws.Cells.Item(171, 4).Value = 44706
$ws.Cells.Item(171, 5).Value = 4
$ws.Cells.Item(171, 6).Value = 100112043
$ws.Cells.Item(171, 7).Value = "Pepino ensalada"
$ws.Cells.Item(171, 8).Value = "Sin especificar"
$ws.Cells.Item(171, 9).Value = "Primera"
$ws.Cells.Item(171, 10).Value = 400
$ws.Cells.Item(171, 11).Value = 15000
$ws.Cells.Item(171, 12).Value = 16000
$ws.Cells.Item(171, 13).Value = 15500
$ws.Cells.Item(171, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(171, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(171, 16).Value = 258
$ws.Cells.Item(171, 17).Value = 60
$ws.Cells.Item(171, 18).Value = "Hortaliza"
